# Insert a new daily price-report row at row 218 (Ajo / Chino / Primera),
# pushing the existing rows 218-287 down to 219-288.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(218).Insert()

$ws.Cells.Item(218, 1).Value  = 9
$ws.Cells.Item(218, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(218, 3).Value  = "Metropolitana"
$ws.Cells.Item(218, 4).Value  = 44875
$ws.Cells.Item(218, 5).Value  = 13
$ws.Cells.Item(218, 6).Value  = 100112003
$ws.Cells.Item(218, 7).Value  = "Ajo"
$ws.Cells.Item(218, 8).Value  = "Chino"
$ws.Cells.Item(218, 9).Value  = "Primera"
$ws.Cells.Item(218, 10).Value = 200
$ws.Cells.Item(218, 11).Value = 14000
$ws.Cells.Item(218, 12).Value = 15000
$ws.Cells.Item(218, 13).Value = 14500
$ws.Cells.Item(218, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(218, 15).Value = "China"
$ws.Cells.Item(218, 16).Value = 1450
$ws.Cells.Item(218, 17).Value = 10
$ws.Cells.Item(218, 18).Value = "Hortaliza"
